$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.159.84"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -7.28%  "
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'3.228.48"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -10.03%  "
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.23%  "
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'174.01"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -16.40%  "
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'508.76"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -10.62%  "
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Value = "'0.584"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -4.59%  "
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').Value = "'  +0.10%  "
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Value = "'3.225.45"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -9.89%  "
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'0.604"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -11.43%  "
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = "'55.85"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -11.71%  "
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Value = "'0.127"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -13.69%  "
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'0.0000249"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -11.24%  "
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'8.93"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -13.77%  "
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'3.754.28"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -9.81%  "
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = "'0.117"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -6.99%  "
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'3.234.61"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -9.87%  "
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'62.984.22"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -7.29%  "
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'16.90"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -11.92%  "
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'10.64"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -12.84%  "
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = "'0.925"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -13.16%  "
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'363.43"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -9.80%  "
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = "'78.35"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -7.41%  "
$ws.Range('E23').Style = 'Normal'

$ws.Range('B24').Value = "'PancakeSwap"
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').Value = "'3.58"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -14.05%  "
$ws.Range('E24').Style = 'Normal'

$ws.Range('B25').Value = "'RenderToken"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'10.71"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -13.23%  "
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = "'5.92"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -3.43%  "
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Value = "'3.70"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -4.18%  "
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = "'2.59"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -10.41%  "
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'11.03"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -11.84%  "
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = "'8.10"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -12.80%  "
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = "'639.80"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -7.05%  "
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = "'27.77"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -11.92%  "
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'6.48"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -14.90%  "
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'10.93"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -10.02%  "
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').Value = "'  -7.40%  "
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Value = "'0.101"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -10.47%  "
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').Value = "'  -0.11%  "
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').Value = "'35.12"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -15.13%  "
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Value = "'0.370"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -9.70%  "
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = "'0.999"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.12%  "
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').Value = "'Maker"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'2.849.09"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -9.83%  "
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = "'Kaspa"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'0.121"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -9.23%  "
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Value = "'0.0₃0635"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -14.99%  "
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'2.60"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -19.31%  "
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').Value = "'2.56"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -7.47%  "
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'2.30"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -13.38%  "
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = "'2.74"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.75%  "
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').Value = "'0.0373"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -9.51%  "
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').Value = "'0.121"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -7.12%  "
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = "'2.87"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -8.44%  "
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = "'130.23"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -6.27%  "
$ws.Range('E51').Style = 'Normal'
